$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: new time-entry ---
$ws.Range("A27").Value = "Lav OC0803 angivSaldoafskrivning"
$ws.Range("B27").Value = "System Analyst "
$ws.Range("C27").Value = "3/3/2020"
$ws.Range("D27").Value = 0.53472222222222221
$ws.Range("E27").Value = 0.60416666666666663
$ws.Range("F27").Value = "1t 30min"

# --- Row 29: new time-entry (string entered before row 28's, to match the
#     original shared-string insertion order) ---
$ws.Range("A29").Value = "Lav SD0802 og DCD0802 angivLineaerAfskrivning"
$ws.Range("B29").Value = "System Analyst "
$ws.Range("C29").Value = "3/5/2020"
$ws.Range("D29").Value = 0.625
$ws.Range("E29").Value = 0.67013888888888884

# --- Row 28: new time-entry ---
$ws.Range("A28").Value = "Review OC0802 og DCD0802"
$ws.Range("B28").Value = "Reviewer"
$ws.Range("C28").Value = "3/4/2020"
$ws.Range("D28").Value = 0.60416666666666663
$ws.Range("E28").Value = 0.625
# H28 loses its running-total formula in the source edit - clear it back out
$ws.Range("H28").ClearContents()

# --- Row 31 gains an (empty) formatted F cell, matching F29/F30's style ---
$ws.Range("F30").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# --- Row 39: extend the shared G formula / H running total down one row ---
$ws.Range("G38:H38").Copy()
$ws.Range("G39:H39").PasteSpecial(-4122)
$ws.Range("G39").Formula = "=E39-D39"
$ws.Range("H39").Formula = "=SUM(G`$5:G39)"

# --- New trailing blank row 54 (mirrors row 53) ---
$ws.Range("C53").Copy()
$ws.Range("C54").PasteSpecial(-4122)

# --- Selection moved to C18 ---
$ws.Range("C18").Select() | Out-Null
